$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(39, 8).Value = 2645908.8
$ws.Cells.Item(39, 9).Value = 3174694.5
$ws.Cells.Item(39, 11).Value = 9524083.5
$ws.Cells.Item(39, 13).Value = -9523787.5

$ws.Cells.Item(99, 8).Value = 1138.6666
$ws.Cells.Item(99, 9).Value = 1253.3
$ws.Cells.Item(99, 10).Value = 565.5
$ws.Cells.Item(99, 11).Value = 3759.9
$ws.Cells.Item(99, 12).Value = 1696.5
$ws.Cells.Item(99, 13).Value = -2261.9
$ws.Cells.Item(99, 14).Value = -4692.5

$ws.Cells.Item(112, 8).Value = 1045
$ws.Cells.Item(112, 9).Value = 759.2857
$ws.Cells.Item(112, 10).Value = 1162.6471
$ws.Cells.Item(112, 11).Value = 2277.8571
$ws.Cells.Item(112, 12).Value = 3487.9413
$ws.Cells.Item(112, 13).Value = -1169.8571
$ws.Cells.Item(112, 14).Value = -5703.9413

$ws.Cells.Item(132, 8).Value = 8628254
$ws.Cells.Item(132, 9).Value = 8936370
$ws.Cells.Item(132, 10).Value = 1000
$ws.Cells.Item(132, 11).Value = 26809110
$ws.Cells.Item(132, 12).Value = 3000
$ws.Cells.Item(132, 13).Value = -26806580
$ws.Cells.Item(132, 14).Value = -8060

$ws.Cells.Item(135, 8).Value = 998.42554
$ws.Cells.Item(135, 9).Value = 627.4737
$ws.Cells.Item(135, 10).Value = 2564.6667
$ws.Cells.Item(135, 11).Value = 5647.263300000001
$ws.Cells.Item(135, 12).Value = 23082.0003
$ws.Cells.Item(135, 13).Value = -3112.263300000001
$ws.Cells.Item(135, 14).Value = -28152.0003

$ws.Cells.Item(137, 8).Value = 1096.0385
$ws.Cells.Item(137, 9).Value = 1118.875
$ws.Cells.Item(137, 10).Value = 1019.9167
$ws.Cells.Item(137, 11).Value = 3356.625
$ws.Cells.Item(137, 12).Value = 3059.7501
$ws.Cells.Item(137, 13).Value = -806.625
$ws.Cells.Item(137, 14).Value = -8159.7501

$ws.Cells.Item(138, 8).Value = 1554.1063
$ws.Cells.Item(138, 9).Value = 1131.5676
$ws.Cells.Item(138, 10).Value = 3117.5
$ws.Cells.Item(138, 11).Value = 3394.7028
$ws.Cells.Item(138, 12).Value = 9352.5
$ws.Cells.Item(138, 13).Value = 1745.2972
$ws.Cells.Item(138, 14).Value = -19632.5

$ws.Cells.Item(141, 8).Value = 1447.0781
$ws.Cells.Item(141, 9).Value = 1273.4138
$ws.Cells.Item(141, 10).Value = 3125.8333
$ws.Cells.Item(141, 11).Value = 3820.2414
$ws.Cells.Item(141, 12).Value = 9377.499899999999
$ws.Cells.Item(141, 13).Value = 1359.7586
$ws.Cells.Item(141, 14).Value = -19737.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1064.77
$ws.Cells.Item(32, 9).Value = 1074.2688
$ws.Cells.Item(32, 11).Value = 1074.2688
$ws.Cells.Item(32, 13).Value = -787.2688000000001

$ws.Cells.Item(44, 8).Value = 12433.3
$ws.Cells.Item(44, 9).Value = 3044
$ws.Cells.Item(44, 10).Value = 13476.556
$ws.Cells.Item(44, 11).Value = 3044
$ws.Cells.Item(44, 12).Value = 13476.556
$ws.Cells.Item(44, 13).Value = -2556
$ws.Cells.Item(44, 14).Value = -14452.556

$ws.Cells.Item(61, 8).Value = 1162.9318
$ws.Cells.Item(61, 9).Value = 1081.0714
$ws.Cells.Item(61, 10).Value = 1306.1875
$ws.Cells.Item(61, 11).Value = 1081.0714
$ws.Cells.Item(61, 12).Value = 1306.1875
$ws.Cells.Item(61, 13).Value = -869.0714
$ws.Cells.Item(61, 14).Value = -1730.1875

$ws.Cells.Item(136, 8).Value = 1162.9318
$ws.Cells.Item(136, 9).Value = 1081.0714
$ws.Cells.Item(136, 10).Value = 1306.1875
$ws.Cells.Item(136, 11).Value = 3243.2142
$ws.Cells.Item(136, 12).Value = 3918.5625
$ws.Cells.Item(136, 13).Value = -693.2142000000003
$ws.Cells.Item(136, 14).Value = -9018.5625

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 2079.509
$ws.Cells.Item(134, 9).Value = 1888.06
$ws.Cells.Item(134, 11).Value = 5664.18
$ws.Cells.Item(134, 13).Value = -3129.18

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 26047.492
$ws.Cells.Item(31, 9).Value = 1375.973
$ws.Cells.Item(31, 10).Value = 67540.5
$ws.Cells.Item(31, 11).Value = 1375.973
$ws.Cells.Item(31, 12).Value = 67540.5
$ws.Cells.Item(31, 13).Value = -1080.973
$ws.Cells.Item(31, 14).Value = -68130.5

$ws.Cells.Item(34, 8).Value = 26047.492
$ws.Cells.Item(34, 9).Value = 1375.973
$ws.Cells.Item(34, 10).Value = 67540.5
$ws.Cells.Item(34, 11).Value = 1375.973
$ws.Cells.Item(34, 12).Value = 67540.5
$ws.Cells.Item(34, 13).Value = -1173.973
$ws.Cells.Item(34, 14).Value = -67944.5

$ws.Cells.Item(88, 8).Value = 30200
$ws.Cells.Item(88, 10).Value = 35266.668
$ws.Cells.Item(88, 12).Value = 35266.668
$ws.Cells.Item(88, 14).Value = -36078.668

$ws.Cells.Item(91, 8).Value = 30200
$ws.Cells.Item(91, 10).Value = 35266.668
$ws.Cells.Item(91, 12).Value = 35266.668
$ws.Cells.Item(91, 14).Value = -38074.668

$ws.Cells.Item(132, 8).Value = 2919.7446
$ws.Cells.Item(132, 9).Value = 2822.3784
$ws.Cells.Item(132, 10).Value = 3280
$ws.Cells.Item(132, 11).Value = 8467.135200000001
$ws.Cells.Item(132, 12).Value = 9840
$ws.Cells.Item(132, 13).Value = -5937.135200000001
$ws.Cells.Item(132, 14).Value = -14900

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 938.70734
$ws.Cells.Item(5, 9).Value = 513.5454999999999
$ws.Cells.Item(5, 10).Value = 1431
$ws.Cells.Item(5, 11).Value = 1540.6365
$ws.Cells.Item(5, 12).Value = 4293
$ws.Cells.Item(5, 13).Value = -1428.6365
$ws.Cells.Item(5, 14).Value = -4517

$ws.Cells.Item(34, 8).Value = 514
$ws.Cells.Item(34, 9).Value = 472.4
$ws.Cells.Item(34, 10).Value = 543.7143
$ws.Cells.Item(34, 11).Value = 1417.2
$ws.Cells.Item(34, 12).Value = 1631.1429
$ws.Cells.Item(34, 13).Value = -1333.2
$ws.Cells.Item(34, 14).Value = -1799.1429

$ws.Cells.Item(56, 8).Value = 4594.143
$ws.Cells.Item(56, 9).Value = 4594.143
$ws.Cells.Item(56, 11).Value = 4594.143
$ws.Cells.Item(56, 13).Value = -4064.143

$ws.Cells.Item(131, 8).Value = 1289.0505
$ws.Cells.Item(131, 10).Value = 1317.4674
$ws.Cells.Item(131, 12).Value = 3952.4022
$ws.Cells.Item(131, 14).Value = -14032.4022

$ws.Cells.Item(135, 8).Value = 938.70734
$ws.Cells.Item(135, 9).Value = 513.5454999999999
$ws.Cells.Item(135, 10).Value = 1431
$ws.Cells.Item(135, 11).Value = 4621.9095
$ws.Cells.Item(135, 12).Value = 12879
$ws.Cells.Item(135, 13).Value = -2086.9095
$ws.Cells.Item(135, 14).Value = -17949

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1081.6757
$ws.Cells.Item(22, 9).Value = 1216.1666
$ws.Cells.Item(22, 10).Value = 1055.6451
$ws.Cells.Item(22, 11).Value = 1216.1666
$ws.Cells.Item(22, 12).Value = 1055.6451
$ws.Cells.Item(22, 13).Value = -921.1666
$ws.Cells.Item(22, 14).Value = -1645.6451

$ws.Cells.Item(27, 8).Value = 1081.6757
$ws.Cells.Item(27, 9).Value = 1216.1666
$ws.Cells.Item(27, 10).Value = 1055.6451
$ws.Cells.Item(27, 11).Value = 1216.1666
$ws.Cells.Item(27, 12).Value = 1055.6451
$ws.Cells.Item(27, 13).Value = -1109.1666
$ws.Cells.Item(27, 14).Value = -1269.6451

$ws.Cells.Item(87, 8).Value = 42000
$ws.Cells.Item(87, 10).Value = 42000
$ws.Cells.Item(87, 12).Value = 42000
$ws.Cells.Item(87, 14).Value = -44246

$ws.Cells.Item(90, 8).Value = 42000
$ws.Cells.Item(90, 10).Value = 42000
$ws.Cells.Item(90, 12).Value = 126000
$ws.Cells.Item(90, 14).Value = -137232

$ws.Cells.Item(132, 8).Value = 2820.111
$ws.Cells.Item(132, 9).Value = 2971.6775
$ws.Cells.Item(132, 10).Value = 1880.4
$ws.Cells.Item(132, 11).Value = 8915.032499999999
$ws.Cells.Item(132, 12).Value = 5641.200000000001
$ws.Cells.Item(132, 13).Value = -6385.032499999999
$ws.Cells.Item(132, 14).Value = -10701.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1590.2307
$ws.Cells.Item(132, 9).Value = 1689.1803
$ws.Cells.Item(132, 10).Value = 1235.1765
$ws.Cells.Item(132, 11).Value = 5067.5409
$ws.Cells.Item(132, 12).Value = 3705.5295
$ws.Cells.Item(132, 13).Value = -2537.5409
$ws.Cells.Item(132, 14).Value = -8765.529500000001

$ws.Cells.Item(136, 8).Value = 503.66666
$ws.Cells.Item(136, 9).Value = 320.4694
$ws.Cells.Item(136, 11).Value = 961.4082000000001
$ws.Cells.Item(136, 13).Value = 1588.5918

Write-Host "Edit applied successfully"
